$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.876.71"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "1.631.30"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'211.71"
$ws.Range("E5").Value = "  -0.07%  "

# Row 6
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("E10").Value = "  -0.80%  "

# Row 11
$ws.Range("E11").Value = "  -0.08%  "

# Row 12
$ws.Range("D12").Value = "1.861.67"
$ws.Range("E12").Value = "  -0.29%  "

# Row 13
$ws.Range("D13").Value = "1.610.71"
$ws.Range("E13").Value = "  -1.51%  "

# Row 14
$ws.Range("E14").Value = "  -0.92%  "

# Row 15
$ws.Range("D15").Value = "'0.557"
$ws.Range("E15").Value = "  -0.79%  "

# Row 16
$ws.Range("D16").Value = "'64.94"
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$ws.Range("D17").Value = "27.908.81"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").Value = "'228.05"
$ws.Range("E18").Value = "  -1.12%  "

# Row 19
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("E20").Value = "  -0.52%  "

# Row 21
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.16%  "

# Row 22
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").Value = "'9.92"
$ws.Range("E23").Value = "  -4.36%  "

# Row 24
$ws.Range("E24").Value = "  -0.40%  "

# Row 25
$ws.Range("D25").Value = "'155.37"
$ws.Range("E25").Value = "  -0.12%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
$ws.Range("D28").Value = "'15.44"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.10%  "

# Row 30
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("E31").Value = "  -0.12%  "

# Row 32
$ws.Range("E32").Value = "  +0.65%  "

# Row 33
$ws.Range("D33").Value = "1.418.74"
$ws.Range("E33").Value = "  +1.47%  "

# Row 34
$ws.Range("E34").Value = "  +1.29%  "

# Row 35
$ws.Range("E35").Value = "  +3.03%  "

# Row 36
$ws.Range("E36").Value = "  -3.66%  "

# Row 37
$ws.Range("E37").Value = "  -1.42%  "

# Row 38
$ws.Range("D38").Value = "'0.0169"
$ws.Range("E38").Value = "  -0.86%  "

# Row 39
$ws.Range("E39").Value = "  -0.23%  "

# Row 40
$ws.Range("D40").Value = "'0.852"
$ws.Range("E40").Value = "  -1.20%  "

# Row 41
$ws.Range("E41").Value = "  -1.70%  "

# Row 42
$ws.Range("D42").Value = "'65.97"
$ws.Range("E42").Value = "  -0.84%  "

# Row 43
$ws.Range("E43").Value = "  -0.57%  "

# Row 44
$ws.Range("E44").Value = "  -0.52%  "

# Row 45
$ws.Range("D45").Value = "1.771.43"
$ws.Range("E45").Value = "  -0.31%  "

# Row 46
$ws.Range("D46").Value = "'2.12"
$ws.Range("E46").Value = "  -3.97%  "

# Row 47
$ws.Range("D47").Value = "'88.65"
$ws.Range("E47").Value = "  +0.80%  "

# Row 48
$ws.Range("E48").Value = "  +1.39%  "

# Row 49
$ws.Range("E49").Value = "  -0.42%  "

# Row 50
$ws.Range("D50").Value = "'7.60"
$ws.Range("E50").Value = "  +1.11%  "

# Row 51
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.16%  "
